$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.325.51'
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").Value = '1.879.77'
$ws.Range("E3").Value = '  -1.69%  '
$ws.Range("E4").Value = '  -0.59%  '
$ws.Range("D5").Value = '246.36'
$ws.Range("E5").Value = '  -3.90%  '
$ws.Range("D6").Value = '0.687'
$ws.Range("E6").Value = '  -6.15%  '
$ws.Range("E7").Value = '  -0.68%  '
$ws.Range("D8").Value = '43.15'
$ws.Range("E8").Value = '  +4.90%  '
$ws.Range("D9").Value = '0.350'
$ws.Range("E9").Value = '  -5.55%  '
$ws.Range("D10").Value = '0.0736'
$ws.Range("E10").Value = '  -3.27%  '
$ws.Range("E11").Value = '  -2.03%  '
$ws.Range("D12").Value = '13.10'
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("D13").Value = '2.149.65'
$ws.Range("E13").Value = '  -1.75%  '
$ws.Range("D14").Value = '0.738'
$ws.Range("E14").Value = '  +0.43%  '
$ws.Range("D15").Value = '4.93'
$ws.Range("E15").Value = '  -1.24%  '
$ws.Range("D16").Value = '1.895.47'
$ws.Range("E16").Value = '  -1.00%  '
$ws.Range("D17").Value = '35.356.77'
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("D18").Value = '73.41'
$ws.Range("E18").Value = '  -2.62%  '
$ws.Range("D19").Value = '0.0₃0821'
$ws.Range("E19").Value = '  -3.22%  '
$ws.Range("D20").Value = '245.44'
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").Value = '12.79'
$ws.Range("E21").Value = '  -2.61%  '
$ws.Range("D22").Value = '4.93'
$ws.Range("E22").Value = '  -4.43%  '
$ws.Range("E23").Value = '  -0.72%  '
$ws.Range("D24").Value = '2.56'
$ws.Range("E24").Value = '  +4.41%  '
$ws.Range("E25").Value = '  -11.38%  '
$ws.Range("D26").Value = '165.40'
$ws.Range("E26").Value = '  -0.53%  '
$ws.Range("D27").Value = '8.45'
$ws.Range("E27").Value = '  -3.64%  '
$ws.Range("D28").Value = '18.27'
$ws.Range("E28").Value = '  -3.13%  '
$ws.Range("E29").Value = '  -4.76%  '
$ws.Range("D30").Value = '4.128.42'
$ws.Range("E31").Value = '  +4.79%  '
$ws.Range("D32").Value = '4.23'
$ws.Range("E32").Value = '  -3.40%  '
$ws.Range("D33").Value = '0.0579'
$ws.Range("E33").Value = '  -2.16%  '
$ws.Range("D34").Value = '4.19'
$ws.Range("E34").Value = '  -1.98%  '
$ws.Range("E35").Value = '  -0.64%  '
$ws.Range("D36").Value = '0.851'
$ws.Range("E36").Value = '  -7.17%  '
$ws.Range("E37").Value = '  -3.53%  '
$ws.Range("D38").Value = '1.57'
$ws.Range("E38").Value = '  -21.50%  '
$ws.Range("D39").Value = '0.0693'
$ws.Range("E39").Value = '  +7.46%  '
$ws.Range("D40").Value = '97.24'
$ws.Range("E40").Value = '  -0.38%  '
$ws.Range("D41").Value = '16.96'
$ws.Range("E41").Value = '  -1.07%  '
$ws.Range("E42").Value = '  -3.09%  '
$ws.Range("D43").Value = '1.08'
$ws.Range("E43").Value = '  -4.35%  '
$ws.Range("D44").Value = '1.288.07'
$ws.Range("E44").Value = '  -4.12%  '
$ws.Range("D45").Value = '2.33'
$ws.Range("E45").Value = '  -6.25%  '
$ws.Range("D46").Value = '0.0811'
$ws.Range("E46").Value = '  +7.27%  '
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("D48").Value = '2.73'
$ws.Range("E48").Value = '  -1.14%  '
$ws.Range("E49").Value = '  +3.07%  '
$ws.Range("D50").Value = '43.15'
$ws.Range("E50").Value = '  -4.82%  '
$ws.Range("E51").Value = '  -7.17%  '
